$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update individual cell values in column A per the diff (rows within 2-201 that changed value)
$ws.Range("A3").Value = 1
$ws.Range("A6").Value = 3
$ws.Range("A7").Value = 1
$ws.Range("A8").Value = 3
$ws.Range("A9").Value = 2
$ws.Range("A10").Value = 3
$ws.Range("A11").Value = 1
$ws.Range("A13").Value = 1
$ws.Range("A14").Value = 3
$ws.Range("A16").Value = 3
$ws.Range("A18").Value = 3
$ws.Range("A19").Value = 3
$ws.Range("A22").Value = 3
$ws.Range("A25").Value = 1
$ws.Range("A26").Value = 2
$ws.Range("A27").Value = 3
$ws.Range("A28").Value = 3
$ws.Range("A29").Value = 3
$ws.Range("A30").Value = 1
$ws.Range("A31").Value = 3
$ws.Range("A33").Value = 3
$ws.Range("A35").Value = 3
$ws.Range("A36").Value = 3
$ws.Range("A37").Value = 3
$ws.Range("A39").Value = 3
$ws.Range("A41").Value = 1
$ws.Range("A42").Value = 2
$ws.Range("A48").Value = 1
$ws.Range("A49").Value = 3
$ws.Range("A50").Value = 3
$ws.Range("A52").Value = 3
$ws.Range("A53").Value = 2
$ws.Range("A54").Value = 2
$ws.Range("A55").Value = 1
$ws.Range("A56").Value = 1
$ws.Range("A57").Value = 1
$ws.Range("A60").Value = 1
$ws.Range("A62").Value = 1
$ws.Range("A66").Value = 3
$ws.Range("A68").Value = 3
$ws.Range("A69").Value = 2
$ws.Range("A70").Value = 1
$ws.Range("A71").Value = 3
$ws.Range("A72").Value = 3
$ws.Range("A73").Value = 2
$ws.Range("A74").Value = 3
$ws.Range("A75").Value = 3
$ws.Range("A76").Value = 3
$ws.Range("A77").Value = 3
$ws.Range("A80").Value = 1
$ws.Range("A81").Value = 3
$ws.Range("A82").Value = 1
$ws.Range("A83").Value = 2
$ws.Range("A86").Value = 3
$ws.Range("A87").Value = 3
$ws.Range("A88").Value = 3
$ws.Range("A92").Value = 1
$ws.Range("A97").Value = 3
$ws.Range("A100").Value = 3
$ws.Range("A102").Value = 2
$ws.Range("A105").Value = 1
$ws.Range("A106").Value = 3
$ws.Range("A107").Value = 1
$ws.Range("A109").Value = 1
$ws.Range("A110").Value = 2
$ws.Range("A111").Value = 1
$ws.Range("A113").Value = 3
$ws.Range("A114").Value = 1
$ws.Range("A115").Value = 1
$ws.Range("A116").Value = 2
$ws.Range("A117").Value = 2
$ws.Range("A118").Value = 3
$ws.Range("A121").Value = 1
$ws.Range("A122").Value = 2
$ws.Range("A123").Value = 3
$ws.Range("A124").Value = 3
$ws.Range("A125").Value = 3
$ws.Range("A126").Value = 2
$ws.Range("A127").Value = 2
$ws.Range("A129").Value = 3
$ws.Range("A130").Value = 2
$ws.Range("A131").Value = 3
$ws.Range("A132").Value = 2
$ws.Range("A133").Value = 3
$ws.Range("A134").Value = 2
$ws.Range("A135").Value = 3
$ws.Range("A136").Value = 2
$ws.Range("A137").Value = 3
$ws.Range("A138").Value = 1
$ws.Range("A141").Value = 3
$ws.Range("A142").Value = 1
$ws.Range("A143").Value = 2
$ws.Range("A146").Value = 2
$ws.Range("A147").Value = 1
$ws.Range("A149").Value = 1
$ws.Range("A150").Value = 2
$ws.Range("A152").Value = 3
$ws.Range("A154").Value = 3
$ws.Range("A156").Value = 3
$ws.Range("A158").Value = 1
$ws.Range("A160").Value = 2
$ws.Range("A162").Value = 3
$ws.Range("A164").Value = 3
$ws.Range("A165").Value = 3
$ws.Range("A167").Value = 1
$ws.Range("A168").Value = 3
$ws.Range("A169").Value = 1
$ws.Range("A170").Value = 3
$ws.Range("A171").Value = 1
$ws.Range("A172").Value = 1
$ws.Range("A174").Value = 3
$ws.Range("A177").Value = 2
$ws.Range("A180").Value = 3
$ws.Range("A181").Value = 2
$ws.Range("A183").Value = 1
$ws.Range("A184").Value = 2
$ws.Range("A186").Value = 3
$ws.Range("A187").Value = 3
$ws.Range("A188").Value = 2
$ws.Range("A189").Value = 3
$ws.Range("A191").Value = 1
$ws.Range("A192").Value = 3
$ws.Range("A194").Value = 3
$ws.Range("A195").Value = 3
$ws.Range("A196").Value = 2
$ws.Range("A197").Value = 3
$ws.Range("A198").Value = 3
$ws.Range("A200").Value = 3
$ws.Range("A201").Value = 3

# Remove the now-unused tail rows (202-251) that existed in the source data
# but are no longer part of the simulation run output.
$ws.Range("A202:A251").ClearContents()
